# Updated symbol list on Sat Jan 28 15:25:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in D (Price) and E (Volume 1h) are stored as text; mark the range as
# Text-formatted before writing so COM keeps the new values as text instead of
# coercing numeric-looking strings (and percentages) into Number cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "306.05"
$ws.Range("E2").Value = "0.42%"
$ws.Range("E3").Value = "7.19%"
$ws.Range("D4").Value = "5.087"
$ws.Range("E4").Value = "1.01%"
$ws.Range("E5").Value = "0.25%"
$ws.Range("D6").Value = "1.946"
$ws.Range("E6").Value = "4.73%"
$ws.Range("D7").Value = "7.934"
$ws.Range("E7").Value = "2.06%"
$ws.Range("D8").Value = "0.9287"
$ws.Range("E8").Value = "0.41%"
$ws.Range("D9").Value = "0.1468"
$ws.Range("E9").Value = "15.79%"
$ws.Range("D10").Value = "0.1930"
$ws.Range("E10").Value = "1.97%"
$ws.Range("D11").Value = "0.09088"
$ws.Range("E11").Value = "0.83%"
$ws.Range("D12").Value = "0.03499"
$ws.Range("E12").Value = "1.56%"
$ws.Range("D13").Value = "0.09792"
$ws.Range("E13").Value = "-0.90%"
$ws.Range("D14").Value = "0.001397"
$ws.Range("E14").Value = "-0.35%"
$ws.Range("D15").Value = "0.005868"
$ws.Range("E15").Value = "-5.95%"
$ws.Range("D16").Value = "3.721"
$ws.Range("E16").Value = "-3.67%"
$ws.Range("D17").Value = "4.184"
$ws.Range("E17").Value = "1.80%"
$ws.Range("D18").Value = "3.467"
$ws.Range("E18").Value = "4.79%"
$ws.Range("D19").Value = "0.3463"
$ws.Range("E19").Value = "1.47%"
$ws.Range("D20").Value = "0.1332"
$ws.Range("E20").Value = "-0.14%"
$ws.Range("D21").Value = "4.808"
$ws.Range("E21").Value = "0.23%"
$ws.Range("D22").Value = "0.2404"
$ws.Range("E22").Value = "2.91%"
$ws.Range("D23").Value = "0.04366"
$ws.Range("E23").Value = "-0.16%"
$ws.Range("E24").Value = "0.02%"
$ws.Range("D25").Value = "0.004266"
$ws.Range("E25").Value = "-12.10%"
$ws.Range("D27").Value = "0.0001299"
$ws.Range("E27").Value = "0.08%"
$ws.Range("D39").Value = "0.02083"
$ws.Range("E39").Value = "6.47%"
$ws.Range("D40").Value = "0.05044"
$ws.Range("E40").Value = "-1.95%"
$ws.Range("D41").Value = "0.007471"
$ws.Range("E41").Value = "-0.54%"
$ws.Range("D42").Value = "0.01009"
$ws.Range("E42").Value = "-0.37%"
$ws.Range("D43").Value = "0.1349"
$ws.Range("E43").Value = "-0.48%"
$ws.Range("D44").Value = "0.002139"
$ws.Range("E44").Value = "1.51%"
$ws.Range("D45").Value = "0.008930"
$ws.Range("E45").Value = "-9.65%"
$ws.Range("D46").Value = "0.00006179"
$ws.Range("E46").Value = "-0.09%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").Value = "0.002803"
$ws.Range("D49").Value = "0.001598"
$ws.Range("E49").Value = "27.81%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.06%"

# Restore default styling so no stray number-format style lingers on the range.
$ws.Range("D2:E51").Style = "Normal"
